{"js": "// Office.js (Word JavaScript API) edit script.\n// This document keeps nearly all of its real content (labels + the\n// vector-calculus formulas) inside floating text boxes / shape groups.\n// The Word JS shape object model in this runtime does not expose a way to\n// reach text nested two levels deep inside a shape group (Group 10 ->\n// Text Box 20), so we perform the edit at the OOXML level: pull the whole\n// body's OOXML, apply the same targeted string surgery the diff describes,\n// and push it back with a full \"replace\" insert. This keeps every other\n// part of the document (styles, other shapes, etc.) untouched while still\n// letting us reach into the nested math zones.\n\nconst body = context.document.body;\nconst ooxmlResult = body.getOoxml();\nawait context.sync();\n\nlet xml = ooxmlResult.value;\n\n// --- Hunk 1: drop the stray _GoBack bookmark at the very start of the body.\nconst staleBookmark = '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>';\nconst staleBookmarkCount = xml.split(staleBookmark).length - 1;\nif (staleBookmarkCount !== 1) {\n  throw new Error(`Expected exactly one leading _GoBack bookmark, found ${staleBookmarkCount}`);\n}\nxml = xml.replace(staleBookmark, \"\");\n\n// --- Hunks 2 & 4: in each of the two \"F\u20d7 = -dU/ds\" math paragraphs, drop the\n// vector-arrow accent around F (leaving a plain \"F\" run) and add a fresh\n// _GoBack bookmark (ids 0 and 1, one per textbox \"story\") right before the\n// math paragraph, matching the bookmark's new home.\nconst vectorArrow = \"\\u20D7\"; // COMBINING RIGHT ARROW ABOVE\nconst accF =\n  '<m:acc><m:accPr><m:chr m:val=\"' + vectorArrow + '\"/><m:ctrlPr><w:rPr>' +\n  '<w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/><w:i/></w:rPr>' +\n  \"</m:ctrlPr></m:accPr><m:e><m:r><w:rPr>\" +\n  '<w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/></w:rPr>' +\n  \"<m:t>F</m:t></m:r></m:e></m:acc>\";\nconst plainF =\n  '<m:r><w:rPr><w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/></w:rPr>' +\n  \"<m:t>F</m:t></m:r>\";\n\nconst targetPrefix = '<w:pPr><w:jc w:val=\"both\"/></w:pPr><m:oMathPara><m:oMath>' + accF;\nconst targetCount = xml.split(targetPrefix).length - 1;\nif (targetCount !== 2) {\n  throw new Error(`Expected exactly two F-vector math paragraphs, found ${targetCount}`);\n}\n\nlet bookmarkId = 0;\nlet rebuilt = \"\";\nlet cursor = 0;\nwhile (true) {\n  const foundAt = xml.indexOf(targetPrefix, cursor);\n  if (foundAt === -1) {\n    rebuilt += xml.slice(cursor);\n    break;\n  }\n  rebuilt += xml.slice(cursor, foundAt);\n  rebuilt +=\n    '<w:pPr><w:jc w:val=\"both\"/></w:pPr>' +\n    `<w:bookmarkStart w:id=\"${bookmarkId}\" w:name=\"_GoBack\"/>` +\n    `<w:bookmarkEnd w:id=\"${bookmarkId}\"/>` +\n    \"<m:oMathPara><m:oMath>\" +\n    plainF;\n  bookmarkId += 1;\n  cursor = foundAt + targetPrefix.length;\n}\nxml = rebuilt;\n\n// --- Hunk 3: add back the (now-missing) v:shapetype definition for the\n// textbox preset used by \"Text Box 19\" 's VML fallback shape.\nconst shapetype =\n  '<v:shapetype id=\"_x0000_t202\" coordsize=\"21600,21600\" o:spt=\"202\" path=\"m,l,21600r21600,l21600,xe\">' +\n  \"<v:stroke joinstyle=\\\"miter\\\"/><v:path gradientshapeok=\\\"t\\\" o:connecttype=\\\"rect\\\"/></v:shapetype>\";\nconst textBox19Marker = '<v:shape id=\"Text Box 19\"';\nconst markerCount = xml.split(textBox19Marker).length - 1;\nif (markerCount !== 1) {\n  throw new Error(`Expected exactly one \"Text Box 19\" VML shape, found ${markerCount}`);\n}\nxml = xml.replace(textBox19Marker, shapetype + textBox19Marker);\n\n// Push the rewritten OOXML back over the whole body.\nbody.insertOoxml(xml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n#\n# Almost all of this document's real content (the labels and the vector-\n# calculus formulas) lives inside floating text boxes that are nested two\n# levels deep inside a shape group (Group 10 -> Text Box 20). The Shapes /\n# GroupItems COM surface in this runtime can see the top-level shapes but\n# cannot descend into that group, so $d.Shapes(...).GroupItems and\n# TextFrame.TextRange don't reach the math zone we need to edit.\n#\n# Instead we round-trip the document through its WordOpenXML: read the\n# single-file XML package, apply the same targeted string surgery the\n# commit's diff describes, and write it back with Range.InsertXML on the\n# whole-document Content range (the documented way to push WordOpenXML\n# back in - see the read-only-property error you get from trying to set\n# Document.WordOpenXML directly).\n\n$d = $word.ActiveDocument\n$xml = $d.Content.WordOpenXML\n\n# --- Hunk 1: drop the stray _GoBack bookmark at the very start of the body.\n$staleBookmark = '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'\n$staleCount = ([regex]::Matches($xml, [regex]::Escape($staleBookmark))).Count\nif ($staleCount -ne 1) {\n    throw (\"Expected exactly one leading _GoBack bookmark, found \" + $staleCount)\n}\n$xml = $xml.Replace($staleBookmark, \"\")\n\n# --- Hunks 2 & 4: in each of the two \"F-with-vector-arrow = -dU/ds\" math\n# paragraphs, drop the vector-arrow accent around F (leaving a plain \"F\"\n# run) and add a fresh _GoBack bookmark (ids 0 and 1, one per textbox\n# \"story\") right before the math paragraph.\n$vectorArrow = [char]0x20D7   # COMBINING RIGHT ARROW ABOVE\n$accF = '<m:acc><m:accPr><m:chr m:val=\"' + $vectorArrow + '\"/><m:ctrlPr><w:rPr>' + `\n    '<w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/><w:i/></w:rPr>' + `\n    '</m:ctrlPr></m:accPr><m:e><m:r><w:rPr>' + `\n    '<w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/></w:rPr>' + `\n    '<m:t>F</m:t></m:r></m:e></m:acc>'\n$plainF = '<m:r><w:rPr><w:rFonts w:ascii=\"Cambria Math\" w:hAnsi=\"Cambria Math\"/></w:rPr>' + `\n    '<m:t>F</m:t></m:r>'\n\n$targetPrefix = '<w:pPr><w:jc w:val=\"both\"/></w:pPr><m:oMathPara><m:oMath>' + $accF\n$targetCount = ([regex]::Matches($xml, [regex]::Escape($targetPrefix))).Count\nif ($targetCount -ne 2) {\n    throw (\"Expected exactly two F-vector math paragraphs, found \" + $targetCount)\n}\n\n$bookmarkId = 0\n$rebuilt = \"\"\n$cursor = 0\nwhile ($true) {\n    $foundAt = $xml.IndexOf($targetPrefix, $cursor)\n    if ($foundAt -eq -1) {\n        $rebuilt += $xml.Substring($cursor)\n        break\n    }\n    $rebuilt += $xml.Substring($cursor, $foundAt - $cursor)\n    $rebuilt += '<w:pPr><w:jc w:val=\"both\"/></w:pPr>' + `\n        ('<w:bookmarkStart w:id=\"' + $bookmarkId + '\" w:name=\"_GoBack\"/>') + `\n        ('<w:bookmarkEnd w:id=\"' + $bookmarkId + '\"/>') + `\n        '<m:oMathPara><m:oMath>' + $plainF\n    $bookmarkId = $bookmarkId + 1\n    $cursor = $foundAt + $targetPrefix.Length\n}\n$xml = $rebuilt\nif ($bookmarkId -ne 2) {\n    throw (\"Expected to insert exactly two bookmarks, inserted \" + $bookmarkId)\n}\n\n# --- Hunk 3: add back the (now-missing) v:shapetype definition for the\n# textbox preset used by \"Text Box 19\" 's VML fallback shape.\n$shapetype = '<v:shapetype id=\"_x0000_t202\" coordsize=\"21600,21600\" o:spt=\"202\" path=\"m,l,21600r21600,l21600,xe\">' + `\n    '<v:stroke joinstyle=\"miter\"/><v:path gradientshapeok=\"t\" o:connecttype=\"rect\"/></v:shapetype>'\n$textBox19Marker = '<v:shape id=\"Text Box 19\"'\n$markerCount = ([regex]::Matches($xml, [regex]::Escape($textBox19Marker))).Count\nif ($markerCount -ne 1) {\n    throw (\"Expected exactly one Text Box 19 VML shape, found \" + $markerCount)\n}\n$xml = $xml.Replace($textBox19Marker, $shapetype + $textBox19Marker)\n\n# Push the rewritten OOXML back over the whole document.\n$d.Content.InsertXML($xml)\n"}
